$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 99, shifting existing rows 99-117 down to 100-118
$ws.Rows.Item(99).Insert()

# Populate the newly inserted row 99 with the new data record
$ws.Range("A99").Value = 11
$ws.Range("B99").Value = 'Vega Monumental Concepción'
$ws.Range("C99").Value = 'Bíobío'
$ws.Range("D99").Value = 44637
$ws.Range("E99").Value = 8
$ws.Range("F99").Value = 'Fruta'
$ws.Range("G99").Value = 100109
$ws.Range("H99").Value = 'Uva'
$ws.Range("I99").Value = 100109001
$ws.Range("J99").Value = 'Uva'
$ws.Range("K99").Value = 'Thompson seedless'
$ws.Range("L99").Value = 'Primera'
$ws.Range("M99").Value = 170
$ws.Range("N99").Value = 9000
$ws.Range("O99").Value = 10000
$ws.Range("P99").Value = 9471
$ws.Range("Q99").Value = '$/bandeja 18 kilos'
$ws.Range("R99").Value = "Región de O'Higgins"
$ws.Range("S99").Value = 526
$ws.Range("T99").Value = 18
